$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1408.3334
$ws.Range("I49").Value = 503.33334
$ws.Range("J49").Value = 2313.3333
$ws.Range("K49").Value = 1510.00002
$ws.Range("L49").Value = 6939.999899999999
$ws.Range("M49").Value = -1374.00002
$ws.Range("N49").Value = -7211.999899999999
$ws.Range("H129").Value = 3084.9333
$ws.Range("I129").Value = 403.33334
$ws.Range("J129").Value = 13811.333
$ws.Range("K129").Value = 1210.00002
$ws.Range("L129").Value = 41433.999
$ws.Range("M129").Value = 3789.99998
$ws.Range("N129").Value = -51433.999
$ws.Range("H137").Value = 27213.309
$ws.Range("I137").Value = 30296.266
$ws.Range("K137").Value = 90888.798
$ws.Range("M137").Value = -88338.798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 46682.668
$ws.Range("J54").Value = 46682.668
$ws.Range("L54").Value = 46682.668
$ws.Range("N54").Value = -48220.668
$ws.Range("H61").Value = 1572.931
$ws.Range("I61").Value = 839.3182
$ws.Range("J61").Value = 3878.5715
$ws.Range("K61").Value = 839.3182
$ws.Range("L61").Value = 3878.5715
$ws.Range("M61").Value = -627.3182
$ws.Range("N61").Value = -4302.5715
$ws.Range("H74").Value = 48923.215
$ws.Range("I74").Value = 92373.95
$ws.Range("J74").Value = 1127.4
$ws.Range("K74").Value = 92373.95
$ws.Range("L74").Value = 1127.4
$ws.Range("M74").Value = -91499.95
$ws.Range("N74").Value = -2875.4
$ws.Range("H77").Value = 48923.215
$ws.Range("I77").Value = 92373.95
$ws.Range("J77").Value = 1127.4
$ws.Range("K77").Value = 461869.75
$ws.Range("L77").Value = 5637
$ws.Range("M77").Value = -457501.75
$ws.Range("N77").Value = -14373
$ws.Range("H132").Value = 1753078
$ws.Range("I132").Value = 2171444.2
$ws.Range("J132").Value = 596418.4
$ws.Range("K132").Value = 6514332.600000001
$ws.Range("L132").Value = 1789255.2
$ws.Range("M132").Value = -6511802.600000001
$ws.Range("N132").Value = -1794315.2
$ws.Range("H136").Value = 1572.931
$ws.Range("I136").Value = 839.3182
$ws.Range("J136").Value = 3878.5715
$ws.Range("K136").Value = 2517.9546
$ws.Range("L136").Value = 11635.7145
$ws.Range("M136").Value = 32.04539999999997
$ws.Range("N136").Value = -16735.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1093.6875
$ws.Range("I22").Value = 9999
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 9999
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -9826
$ws.Range("N22").Value = -846
$ws.Range("H134").Value = 20026.637
$ws.Range("I134").Value = 1138.6666
$ws.Range("J134").Value = 81049.30499999999
$ws.Range("K134").Value = 3415.9998
$ws.Range("L134").Value = 243147.915
$ws.Range("M134").Value = -880.9998000000001
$ws.Range("N134").Value = -248217.915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30278.361
$ws.Range("I31").Value = 33132.934
$ws.Range("J31").Value = 12580
$ws.Range("K31").Value = 33132.934
$ws.Range("L31").Value = 12580
$ws.Range("M31").Value = -32837.934
$ws.Range("N31").Value = -13170
$ws.Range("H34").Value = 30278.361
$ws.Range("I34").Value = 33132.934
$ws.Range("J34").Value = 12580
$ws.Range("K34").Value = 33132.934
$ws.Range("L34").Value = 12580
$ws.Range("M34").Value = -32930.934
$ws.Range("N34").Value = -12984
$ws.Range("H58").Value = 4000.7
$ws.Range("I58").Value = 1226.2
$ws.Range("J58").Value = 6775.2
$ws.Range("K58").Value = 1226.2
$ws.Range("L58").Value = 6775.2
$ws.Range("M58").Value = -1023.2
$ws.Range("N58").Value = -7181.2
$ws.Range("H132").Value = 2276.75
$ws.Range("I132").Value = 1654.32
$ws.Range("J132").Value = 4499.7144
$ws.Range("K132").Value = 4962.96
$ws.Range("L132").Value = 13499.1432
$ws.Range("M132").Value = -2432.96
$ws.Range("N132").Value = -18559.1432
$ws.Range("H134").Value = 13159565
$ws.Range("I134").Value = 1266.4
$ws.Range("J134").Value = 38463984
$ws.Range("K134").Value = 3799.2
$ws.Range("L134").Value = 115391952
$ws.Range("M134").Value = -1264.2
$ws.Range("N134").Value = -115397022
$ws.Range("H136").Value = 4000.7
$ws.Range("I136").Value = 1226.2
$ws.Range("J136").Value = 6775.2
$ws.Range("K136").Value = 3678.6
$ws.Range("L136").Value = 20325.6
$ws.Range("M136").Value = -1128.6
$ws.Range("N136").Value = -25425.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 166686320
$ws.Range("J39").Value = 250029010
$ws.Range("L39").Value = 750087030
$ws.Range("N39").Value = -750087618
$ws.Range("H62").Value = 3376.625
$ws.Range("I62").Value = 450
$ws.Range("J62").Value = 3794.7144
$ws.Range("K62").Value = 1350
$ws.Range("L62").Value = 11384.1432
$ws.Range("M62").Value = -664
$ws.Range("N62").Value = -12756.1432
$ws.Range("H64").Value = 2168525.5
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 2529529.8
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 7588589.399999999
$ws.Range("M64").Value = -7230
$ws.Range("N64").Value = -7589129.399999999
$ws.Range("H65").Value = 3376.625
$ws.Range("I65").Value = 450
$ws.Range("J65").Value = 3794.7144
$ws.Range("K65").Value = 4050
$ws.Range("L65").Value = 34152.4296
$ws.Range("M65").Value = -618
$ws.Range("N65").Value = -41016.4296
$ws.Range("H67").Value = 2168525.5
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 2529529.8
$ws.Range("K67").Value = 7500
$ws.Range("L67").Value = 7588589.399999999
$ws.Range("M67").Value = -6564
$ws.Range("N67").Value = -7590461.399999999
$ws.Range("H68").Value = 12105.444
$ws.Range("I68").Value = 50389.5
$ws.Range("J68").Value = 1167.1428
$ws.Range("K68").Value = 151168.5
$ws.Range("L68").Value = 3501.4284
$ws.Range("M68").Value = -150357.5
$ws.Range("N68").Value = -5123.428400000001
$ws.Range("H71").Value = 12105.444
$ws.Range("I71").Value = 50389.5
$ws.Range("J71").Value = 1167.1428
$ws.Range("K71").Value = 453505.5
$ws.Range("L71").Value = 10504.2852
$ws.Range("M71").Value = -449449.5
$ws.Range("N71").Value = -18616.2852
$ws.Range("H107").Value = 232.9
$ws.Range("I107").Value = 122
$ws.Range("J107").Value = 245.22223
$ws.Range("K107").Value = 366
$ws.Range("L107").Value = 735.66669
$ws.Range("M107").Value = 1554
$ws.Range("N107").Value = -4575.66669
$ws.Range("H127").Value = 2565522.5
$ws.Range("J127").Value = 2565522.5
$ws.Range("L127").Value = 7696567.5
$ws.Range("N127").Value = -7706487.5
$ws.Range("H131").Value = 30449528
$ws.Range("J131").Value = 34421130
$ws.Range("L131").Value = 103263390
$ws.Range("N131").Value = -103273470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4210.552
$ws.Range("I122").Value = 4035.2942
$ws.Range("K122").Value = 12105.8826
$ws.Range("M122").Value = -9655.882599999999
$ws.Range("H132").Value = 40015.27
$ws.Range("I132").Value = 1222.0454
$ws.Range("J132").Value = 253378
$ws.Range("K132").Value = 3666.1362
$ws.Range("L132").Value = 760134
$ws.Range("M132").Value = -1136.1362
$ws.Range("N132").Value = -765194

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 288506
$ws.Range("I132").Value = 73592.78999999999
$ws.Range("J132").Value = 718332.4399999999
$ws.Range("K132").Value = 220778.37
$ws.Range("L132").Value = 2154997.32
$ws.Range("M132").Value = -218248.37
$ws.Range("N132").Value = -2160057.32
$ws.Range("H136").Value = 590501
$ws.Range("I136").Value = 1112488.5
$ws.Range("J136").Value = 3265
$ws.Range("K136").Value = 3337465.5
$ws.Range("L136").Value = 9795
$ws.Range("M136").Value = -3334915.5
$ws.Range("N136").Value = -14895

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 36666.668
$ws.Range("J58").Value = 36666.668
$ws.Range("L58").Value = 36666.668
$ws.Range("N58").Value = -37282.668
$ws.Range("H132").Value = 6167.6816
$ws.Range("I132").Value = 1432.4706
$ws.Range("J132").Value = 22267.4
$ws.Range("K132").Value = 4297.4118
$ws.Range("L132").Value = 66802.20000000001
$ws.Range("M132").Value = -1767.4118
$ws.Range("N132").Value = -71862.20000000001
$ws.Range("H136").Value = 4526030
$ws.Range("I136").Value = 4204031.5
$ws.Range("J136").Value = 10000005
$ws.Range("K136").Value = 12612094.5
$ws.Range("L136").Value = 30000015
$ws.Range("M136").Value = -12609544.5
$ws.Range("N136").Value = -30005115
